$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.830.98'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '2.384.54'
$ws.Range("E3").Value = '  -0.69%  '
$__style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = $__style
$ws.Range("E4").Value = '  -0.13%  '
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '507.51'
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = '  +0.48%  '
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.36'
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -1.18%  '
$ws.Range("D9").Value = '2.398.23'
$ws.Range("E9").Value = '  -0.83%  '
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0990'
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = '  +2.84%  '
$ws.Range("E11").Value = '  +0.34%  '
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.85'
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = '  +5.77%  '
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.328'
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = '  +2.25%  '
$ws.Range("D14").Value = '2.807.67'
$ws.Range("E14").Value = '  -0.83%  '
$ws.Range("D15").Value = '56.701.13'
$ws.Range("E15").Value = '  -0.30%  '
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.72'
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("E17").Value = '  +0.32%  '
$ws.Range("D18").Value = '2.405.59'
$ws.Range("E18").Value = '  +2.37%  '
$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.14'
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = '  -0.36%  '
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.08'
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = '  +1.26%  '
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '311.17'
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("E22").Value = '  -0.33%  '
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = '  +0.03%  '
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.64'
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = '  +2.07%  '
$ws.Range("E25").Value = '  +0.04%  '
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.372'
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = '  -1.69%  '
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.148'
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = '  -3.75%  '
$ws.Range("E28").Value = '  -2.25%  '
$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.35'
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = '  +1.70%  '
$ws.Range("E30").Value = '  -0.20%  '
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.66'
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = '  -0.77%  '
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.87'
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = '  -0.78%  '
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.08'
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = '  -2.78%  '
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("E38").Value = '  -1.70%  '
$ws.Range("E39").Value = '  +4.48%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("E41").Value = '  -2.29%  '
$ws.Range("E42").Value = '  +1.71%  '
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '129.39'
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("E45").Value = '  +0.88%  '
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0900'
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = '  -0.89%  '
$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '241.78'
$ws.Range("D47").Style = $__style
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0485'
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = '  -0.53%  '
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0208'
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("E51").Value = '  -0.90%  '
